$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 9029.763000000001
$ws.Cells.Item(137, 9).Value = 14292.048
$ws.Cells.Item(137, 10).Value = 2529.2942
$ws.Cells.Item(137, 11).Value = 42876.144
$ws.Cells.Item(137, 12).Value = 7587.882599999999
$ws.Cells.Item(137, 13).Value = -40326.144
$ws.Cells.Item(137, 14).Value = -12687.8826

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6574.7534
$ws.Cells.Item(32, 9).Value = 6285.1
$ws.Cells.Item(32, 11).Value = 6285.1
$ws.Cells.Item(32, 13).Value = -5998.1

$ws.Cells.Item(102, 8).Value = 15138.207
$ws.Cells.Item(102, 9).Value = 21325.938
$ws.Cells.Item(102, 10).Value = 7522.5386
$ws.Cells.Item(102, 11).Value = 21325.938
$ws.Cells.Item(102, 12).Value = 7522.5386
$ws.Cells.Item(102, 13).Value = -19703.938
$ws.Cells.Item(102, 14).Value = -10766.5386

$ws.Cells.Item(122, 8).Value = 1116989.1
$ws.Cells.Item(122, 9).Value = 5717.1113
$ws.Cells.Item(122, 11).Value = 17151.3339
$ws.Cells.Item(122, 13).Value = -14701.3339

$ws.Cells.Item(132, 8).Value = 2789.4546
$ws.Cells.Item(132, 9).Value = 1844.7097
$ws.Cells.Item(132, 11).Value = 5534.1291
$ws.Cells.Item(132, 13).Value = -3004.1291

$ws.Cells.Item(135, 8).Value = 70121.375
$ws.Cells.Item(135, 10).Value = 70121.375
$ws.Cells.Item(135, 12).Value = 70121.375
$ws.Cells.Item(135, 14).Value = -80261.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 12065.186
$ws.Cells.Item(99, 9).Value = 15344.389
$ws.Cells.Item(99, 10).Value = 5506.778
$ws.Cells.Item(99, 11).Value = 15344.389
$ws.Cells.Item(99, 12).Value = 5506.778
$ws.Cells.Item(99, 13).Value = -13846.389
$ws.Cells.Item(99, 14).Value = -8502.778

$ws.Cells.Item(107, 8).Value = 2631.3845
$ws.Cells.Item(107, 9).Value = 2684
$ws.Cells.Item(107, 11).Value = 2684
$ws.Cells.Item(107, 13).Value = -764

$ws.Cells.Item(130, 8).Value = 80787.69500000001
$ws.Cells.Item(130, 10).Value = 80787.69500000001
$ws.Cells.Item(130, 12).Value = 80787.69500000001
$ws.Cells.Item(130, 14).Value = -90827.69500000001

$ws.Cells.Item(134, 8).Value = 6742.4116
$ws.Cells.Item(134, 9).Value = 6884.448
$ws.Cells.Item(134, 10).Value = 5918.6
$ws.Cells.Item(134, 11).Value = 20653.344
$ws.Cells.Item(134, 12).Value = 17755.8
$ws.Cells.Item(134, 13).Value = -18118.344
$ws.Cells.Item(134, 14).Value = -22825.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 180.5
$ws.Cells.Item(7, 9).Value = 180.5
$ws.Cells.Item(7, 11).Value = 180.5
$ws.Cells.Item(7, 13).Value = -67.5

$ws.Cells.Item(16, 8).Value = 53662.21
$ws.Cells.Item(16, 9).Value = 1042.5454
$ws.Cells.Item(16, 10).Value = 126014.25
$ws.Cells.Item(16, 11).Value = 1042.5454
$ws.Cells.Item(16, 12).Value = 126014.25
$ws.Cells.Item(16, 13).Value = -755.5454
$ws.Cells.Item(16, 14).Value = -126588.25

$ws.Cells.Item(31, 8).Value = 6152.1953
$ws.Cells.Item(31, 9).Value = 6367.75
$ws.Cells.Item(31, 10).Value = 5385.778
$ws.Cells.Item(31, 11).Value = 6367.75
$ws.Cells.Item(31, 12).Value = 5385.778
$ws.Cells.Item(31, 13).Value = -6072.75
$ws.Cells.Item(31, 14).Value = -5975.778

$ws.Cells.Item(34, 8).Value = 6152.1953
$ws.Cells.Item(34, 9).Value = 6367.75
$ws.Cells.Item(34, 10).Value = 5385.778
$ws.Cells.Item(34, 11).Value = 6367.75
$ws.Cells.Item(34, 12).Value = 5385.778
$ws.Cells.Item(34, 13).Value = -6165.75
$ws.Cells.Item(34, 14).Value = -5789.778

$ws.Cells.Item(58, 8).Value = 2727.8438
$ws.Cells.Item(58, 9).Value = 3206.375
$ws.Cells.Item(58, 10).Value = 2249.3125
$ws.Cells.Item(58, 11).Value = 3206.375
$ws.Cells.Item(58, 12).Value = 2249.3125
$ws.Cells.Item(58, 13).Value = -3003.375
$ws.Cells.Item(58, 14).Value = -2655.3125

$ws.Cells.Item(105, 8).Value = 124595.586
$ws.Cells.Item(105, 10).Value = 1093.25
$ws.Cells.Item(105, 12).Value = 1093.25
$ws.Cells.Item(105, 14).Value = -4587.25

$ws.Cells.Item(113, 8).Value = 53662.21
$ws.Cells.Item(113, 9).Value = 1042.5454
$ws.Cells.Item(113, 10).Value = 126014.25
$ws.Cells.Item(113, 11).Value = 1042.5454
$ws.Cells.Item(113, 12).Value = 126014.25
$ws.Cells.Item(113, 13).Value = 1127.4546
$ws.Cells.Item(113, 14).Value = -130354.25

$ws.Cells.Item(122, 8).Value = 9787.286
$ws.Cells.Item(122, 9).Value = 11897.546
$ws.Cells.Item(122, 11).Value = 35692.638
$ws.Cells.Item(122, 13).Value = -33242.638

$ws.Cells.Item(132, 8).Value = 1729.1904
$ws.Cells.Item(132, 9).Value = 1508.4
$ws.Cells.Item(132, 10).Value = 2281.1667
$ws.Cells.Item(132, 11).Value = 4525.200000000001
$ws.Cells.Item(132, 12).Value = 6843.500100000001
$ws.Cells.Item(132, 13).Value = -1995.200000000001
$ws.Cells.Item(132, 14).Value = -11903.5001

$ws.Cells.Item(134, 8).Value = 4731.613
$ws.Cells.Item(134, 9).Value = 5472.577
$ws.Cells.Item(134, 11).Value = 16417.731
$ws.Cells.Item(134, 13).Value = -13882.731

$ws.Cells.Item(136, 8).Value = 2727.8438
$ws.Cells.Item(136, 9).Value = 3206.375
$ws.Cells.Item(136, 10).Value = 2249.3125
$ws.Cells.Item(136, 11).Value = 9619.125
$ws.Cells.Item(136, 12).Value = 6747.9375
$ws.Cells.Item(136, 13).Value = -7069.125
$ws.Cells.Item(136, 14).Value = -11847.9375

$ws.Cells.Item(141, 8).Value = 247671.61
$ws.Cells.Item(141, 10).Value = 262060.92
$ws.Cells.Item(141, 12).Value = 262060.92
$ws.Cells.Item(141, 14).Value = -272420.92

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(56, 8).Value = 6647.7144
$ws.Cells.Item(56, 9).Value = 6647.7144
$ws.Cells.Item(56, 11).Value = 6647.7144
$ws.Cells.Item(56, 13).Value = -6117.7144

$ws.Cells.Item(92, 8).Value = 750
$ws.Cells.Item(92, 10).Value = 750
$ws.Cells.Item(92, 12).Value = 2250
$ws.Cells.Item(92, 14).Value = -4746

$ws.Cells.Item(97, 8).Value = 120101
$ws.Cells.Item(97, 9).Value = 120101
$ws.Cells.Item(97, 11).Value = 360303
$ws.Cells.Item(97, 13).Value = -359807

$ws.Cells.Item(140, 8).Value = 1503.6154
$ws.Cells.Item(140, 9).Value = 1503.6154
$ws.Cells.Item(140, 11).Value = 4510.8462
$ws.Cells.Item(140, 13).Value = 669.1538

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 2233.3333

$ws.Cells.Item(113, 8).Value = 10293.923
$ws.Cells.Item(113, 9).Value = 12282.1
$ws.Cells.Item(113, 10).Value = 3666.6667
$ws.Cells.Item(113, 11).Value = 12282.1
$ws.Cells.Item(113, 12).Value = 3666.6667
$ws.Cells.Item(113, 13).Value = -10112.1
$ws.Cells.Item(113, 14).Value = -8006.6667

$ws.Cells.Item(132, 8).Value = 4046.7693
$ws.Cells.Item(132, 9).Value = 4250.1304
$ws.Cells.Item(132, 10).Value = 2487.6667
$ws.Cells.Item(132, 11).Value = 12750.3912
$ws.Cells.Item(132, 12).Value = 7463.000100000001
$ws.Cells.Item(132, 13).Value = -10220.3912
$ws.Cells.Item(132, 14).Value = -12523.0001

$ws.Cells.Item(138, 8).Value = 97999
$ws.Cells.Item(138, 10).Value = 97999
$ws.Cells.Item(138, 12).Value = 97999
$ws.Cells.Item(138, 14).Value = -108279

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 598540.9399999999
$ws.Cells.Item(132, 9).Value = 1065405.1
$ws.Cells.Item(132, 10).Value = 4350.273
$ws.Cells.Item(132, 11).Value = 3196215.3
$ws.Cells.Item(132, 12).Value = 13050.819
$ws.Cells.Item(132, 13).Value = -3193685.3
$ws.Cells.Item(132, 14).Value = -18110.819

$ws.Cells.Item(136, 8).Value = 5687.0454
$ws.Cells.Item(136, 9).Value = 4739.6
$ws.Cells.Item(136, 10).Value = 5965.706
$ws.Cells.Item(136, 11).Value = 14218.8
$ws.Cells.Item(136, 12).Value = 17897.118
$ws.Cells.Item(136, 13).Value = -11668.8
$ws.Cells.Item(136, 14).Value = -22997.118

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 7780.2344
$ws.Cells.Item(132, 9).Value = 9098.325000000001
$ws.Cells.Item(132, 11).Value = 27294.975
$ws.Cells.Item(132, 13).Value = -24764.975

$ws.Cells.Item(136, 8).Value = 558698.6
$ws.Cells.Item(136, 9).Value = 776728.4399999999
$ws.Cells.Item(136, 10).Value = 13624.125
$ws.Cells.Item(136, 11).Value = 2330185.32
$ws.Cells.Item(136, 12).Value = 40872.375
$ws.Cells.Item(136, 13).Value = -2327635.32
$ws.Cells.Item(136, 14).Value = -45972.375

$ws.Cells.Item(141, 8).Value = 92124.25
$ws.Cells.Item(141, 10).Value = 83999.14
$ws.Cells.Item(141, 12).Value = 83999.14
$ws.Cells.Item(141, 14).Value = -94359.14
